$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = 'ISIN'
$ws.Cells.Item(1,2).Value = 'Stock Name'
$ws.Cells.Item(1,3).Value = 'Mutual Fund'
$ws.Cells.Item(1,4).Value = 'Jan_2026'
$ws.Cells.Item(1,5).Value = 'Dec_2025'
$ws.Cells.Item(1,6).Value = 'Nov_2025'
$ws.Cells.Item(1,7).Value = 'MoM'
$ws.Cells.Item(1,8).Value = 'QoQ'

$ws.Cells.Item(2,1).Value = 'INE040A16GJ4'
$ws.Cells.Item(2,2).Value = 'HDFC Bank Ltd CD 25-Feb-2026'
$ws.Cells.Item(2,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(2,4).Value = 7.89811
$ws.Cells.Item(2,5).Value = 8.269011
$ws.Cells.Item(2,6).Value = 5.771966
$ws.Cells.Item(2,7).Value = -0.3709010000000008
$ws.Cells.Item(2,8).Value = 2.126144

$ws.Cells.Item(3,1).Value = 'INE237A166Z3'
$ws.Cells.Item(3,2).Value = 'Kotak Mahindra Bank Ltd CD 27-Feb-2026'
$ws.Cells.Item(3,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(3,4).Value = 7.895722
$ws.Cells.Item(3,5).Value = 8.266792
$ws.Cells.Item(3,6).Value = 0
$ws.Cells.Item(3,7).Value = -0.3710700000000005
$ws.Cells.Item(3,8).Value = 7.895722

$ws.Cells.Item(4,1).Value = 'INE028A16HW0'
$ws.Cells.Item(4,2).Value = 'Bank Of Baroda CD 06-Mar-2026'
$ws.Cells.Item(4,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(4,4).Value = 7.884695
$ws.Cells.Item(4,5).Value = 8.257506
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = -0.3728109999999996
$ws.Cells.Item(4,8).Value = 7.884695

$ws.Cells.Item(5,1).Value = 'INE508G14IE9'
$ws.Cells.Item(5,2).Value = 'Time Technoplast Limited CP 20-Mar-2026'
$ws.Cells.Item(5,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(5,4).Value = 7.8484
$ws.Cells.Item(5,5).Value = 8.211885
$ws.Cells.Item(5,6).Value = 0
$ws.Cells.Item(5,7).Value = -0.3634850000000007
$ws.Cells.Item(5,8).Value = 7.8484

$ws.Cells.Item(6,1).Value = 'INE296A14E79'
$ws.Cells.Item(6,2).Value = 'Bajaj Finance Limited CP 29-Apr-2026'
$ws.Cells.Item(6,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(6,4).Value = 7.787876
$ws.Cells.Item(6,5).Value = 0
$ws.Cells.Item(6,6).Value = 0
$ws.Cells.Item(6,7).Value = 7.787876
$ws.Cells.Item(6,8).Value = 7.787876

$ws.Cells.Item(7,1).Value = 'INE763G14E51'
$ws.Cells.Item(7,2).Value = 'ICICI Securities Ltd CP 30-Apr-2026'
$ws.Cells.Item(7,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(7,4).Value = 7.782529
$ws.Cells.Item(7,5).Value = 0
$ws.Cells.Item(7,6).Value = 0
$ws.Cells.Item(7,7).Value = 7.782529
$ws.Cells.Item(7,8).Value = 7.782529

$ws.Cells.Item(8,1).Value = 'INE233A146K1'
$ws.Cells.Item(8,2).Value = 'Godrej Industries Ltd CP 20-Apr-2026'
$ws.Cells.Item(8,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(8,4).Value = 5.854878
$ws.Cells.Item(8,5).Value = 0
$ws.Cells.Item(8,6).Value = 0
$ws.Cells.Item(8,7).Value = 5.854878
$ws.Cells.Item(8,8).Value = 5.854878

$ws.Cells.Item(9,1).Value = 'INE476A16A24'
$ws.Cells.Item(9,2).Value = 'Canara Bank CD 03-Feb-2026'
$ws.Cells.Item(9,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(9,4).Value = 3.964613
$ws.Cells.Item(9,5).Value = 4.149348
$ws.Cells.Item(9,6).Value = 2.89556
$ws.Cells.Item(9,7).Value = -0.1847349999999999
$ws.Cells.Item(9,8).Value = 1.069053

$ws.Cells.Item(10,1).Value = 'INE514E16CJ9'
$ws.Cells.Item(10,2).Value = 'EXIM Bank CD 04-Mar-2026'
$ws.Cells.Item(10,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(10,4).Value = 3.943704
$ws.Cells.Item(10,5).Value = 4.130146
$ws.Cells.Item(10,6).Value = 0
$ws.Cells.Item(10,7).Value = -0.186442
$ws.Cells.Item(10,8).Value = 3.943704

$ws.Cells.Item(11,1).Value = 'INE115A14FK9'
$ws.Cells.Item(11,2).Value = 'LIC Housing Finance Ltd CP 11-Mar-2026'
$ws.Cells.Item(11,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(11,4).Value = 3.93827
$ws.Cells.Item(11,5).Value = 4.12479
$ws.Cells.Item(11,6).Value = 0
$ws.Cells.Item(11,7).Value = -0.1865199999999998
$ws.Cells.Item(11,8).Value = 3.93827

$ws.Cells.Item(12,1).Value = 'INE028E14TY8'
$ws.Cells.Item(12,2).Value = 'Kotak Securities Ltd CP 10-Mar-2026'
$ws.Cells.Item(12,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(12,4).Value = 3.937286
$ws.Cells.Item(12,5).Value = 4.121816
$ws.Cells.Item(12,6).Value = 0
$ws.Cells.Item(12,7).Value = -0.1845300000000001
$ws.Cells.Item(12,8).Value = 3.937286

$ws.Cells.Item(13,1).Value = 'INE134E14AX6'
$ws.Cells.Item(13,2).Value = 'Power Finance Corp Ltd CP 15-Apr-2026'
$ws.Cells.Item(13,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(13,4).Value = 3.909039
$ws.Cells.Item(13,5).Value = 0
$ws.Cells.Item(13,6).Value = 0
$ws.Cells.Item(13,7).Value = 3.909039
$ws.Cells.Item(13,8).Value = 3.909039

$ws.Cells.Item(14,1).Value = 'INE238AD6BX7'
$ws.Cells.Item(14,2).Value = 'Axis Bank Limited CD 29-Apr-2026'
$ws.Cells.Item(14,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(14,4).Value = 3.898289
$ws.Cells.Item(14,5).Value = 0
$ws.Cells.Item(14,6).Value = 0
$ws.Cells.Item(14,7).Value = 3.898289
$ws.Cells.Item(14,8).Value = 3.898289

$ws.Cells.Item(15,1).Value = 'INE556F16AX2'
$ws.Cells.Item(15,2).Value = 'SIDBI CD 05-Dec-2025'
$ws.Cells.Item(15,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 0
$ws.Cells.Item(15,6).Value = 5.848539
$ws.Cells.Item(15,7).Value = 0
$ws.Cells.Item(15,8).Value = -5.848539

$ws.Cells.Item(16,1).Value = 'INE414G14UB1'
$ws.Cells.Item(16,2).Value = 'Muthoot Finance Ltd CP 12-Jan-2026'
$ws.Cells.Item(16,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(16,4).Value = 0
$ws.Cells.Item(16,5).Value = 4.163189
$ws.Cells.Item(16,6).Value = 2.903499
$ws.Cells.Item(16,7).Value = -4.163189
$ws.Cells.Item(16,8).Value = -2.903499

$ws.Cells.Item(17,1).Value = 'INE618R14018'
$ws.Cells.Item(17,2).Value = 'SG Finserv Ltd CP 27-Jan-2026'
$ws.Cells.Item(17,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(17,4).Value = 0
$ws.Cells.Item(17,5).Value = 4.149965
$ws.Cells.Item(17,6).Value = 2.892022
$ws.Cells.Item(17,7).Value = -4.149965
$ws.Cells.Item(17,8).Value = -2.892022

$ws.Cells.Item(18,1).Value = 'INE763G14B62'
$ws.Cells.Item(18,2).Value = 'ICICI Securities Ltd CP 12-Jan-2026'
$ws.Cells.Item(18,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(18,4).Value = 0
$ws.Cells.Item(18,5).Value = 0
$ws.Cells.Item(18,6).Value = 5.807624
$ws.Cells.Item(18,7).Value = 0
$ws.Cells.Item(18,8).Value = -5.807624

$ws.Cells.Item(19,1).Value = 'INE508G14HS1'
$ws.Cells.Item(19,2).Value = 'Time Technoplast Limited CP 24-Dec-2025'
$ws.Cells.Item(19,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(19,4).Value = 0
$ws.Cells.Item(19,5).Value = 0
$ws.Cells.Item(19,6).Value = 5.824888
$ws.Cells.Item(19,7).Value = 0
$ws.Cells.Item(19,8).Value = -5.824888

$ws.Cells.Item(20,1).Value = 'INE018A14LR8'
$ws.Cells.Item(20,2).Value = 'Larsen & Toubro Ltd CP 26-Dec-2025'
$ws.Cells.Item(20,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(20,4).Value = 0
$ws.Cells.Item(20,5).Value = 0
$ws.Cells.Item(20,6).Value = 5.828793
$ws.Cells.Item(20,7).Value = 0
$ws.Cells.Item(20,8).Value = -5.828793

$ws.Cells.Item(21,1).Value = 'INE296A14C71'
$ws.Cells.Item(21,2).Value = 'Bajaj Finance Limited CP 13-Jan-2026'
$ws.Cells.Item(21,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(21,4).Value = 0
$ws.Cells.Item(21,5).Value = 0
$ws.Cells.Item(21,6).Value = 2.903587
$ws.Cells.Item(21,7).Value = 0
$ws.Cells.Item(21,8).Value = -2.903587

$ws.Cells.Item(22,1).Value = 'INE261F16900'
$ws.Cells.Item(22,2).Value = 'NABARD CD 22-Jan-2026'
$ws.Cells.Item(22,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(22,4).Value = 0
$ws.Cells.Item(22,5).Value = 0
$ws.Cells.Item(22,6).Value = 2.90097
$ws.Cells.Item(22,7).Value = 0
$ws.Cells.Item(22,8).Value = -2.90097

$ws.Cells.Item(23,1).Value = 'INE238AD6AE9'
$ws.Cells.Item(23,2).Value = 'Axis Bank Limited CD 08-Jan-2026'
$ws.Cells.Item(23,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(23,4).Value = 0
$ws.Cells.Item(23,5).Value = 4.166773
$ws.Cells.Item(23,6).Value = 2.907702
$ws.Cells.Item(23,7).Value = -4.166773
$ws.Cells.Item(23,8).Value = -2.907702

$ws.Cells.Item(24,1).Value = 'INE233A144Q3'
$ws.Cells.Item(24,2).Value = 'Godrej Industries Ltd CP 16-Jan-2026'
$ws.Cells.Item(24,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(24,4).Value = 0
$ws.Cells.Item(24,5).Value = 6.241449
$ws.Cells.Item(24,6).Value = 4.354897
$ws.Cells.Item(24,7).Value = -6.241449
$ws.Cells.Item(24,8).Value = -4.354897

$ws.Cells.Item(25,1).Value = 'INE192R14287'
$ws.Cells.Item(25,2).Value = 'Avenue Supermarts CP 29-Dec-2025'
$ws.Cells.Item(25,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(25,5).Value = 0
$ws.Cells.Item(25,6).Value = 5.825315
$ws.Cells.Item(25,7).Value = 0
$ws.Cells.Item(25,8).Value = -5.825315

$ws.Cells.Item(26,1).Value = 'INE033L14OJ5'
$ws.Cells.Item(26,2).Value = 'Tata Capital Housing Fin CP 14-Jan-2026'
$ws.Cells.Item(26,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(26,4).Value = 0
$ws.Cells.Item(26,5).Value = 8.325077
$ws.Cells.Item(26,6).Value = 5.809076
$ws.Cells.Item(26,7).Value = -8.325077
$ws.Cells.Item(26,8).Value = -5.809076

$ws.Cells.Item(27,1).Value = 'INE028A16JS4'
$ws.Cells.Item(27,2).Value = 'Bank Of Baroda CD 02-Dec-2025'
$ws.Cells.Item(27,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(27,4).Value = 0
$ws.Cells.Item(27,5).Value = 0
$ws.Cells.Item(27,6).Value = 2.925713
$ws.Cells.Item(27,7).Value = 0
$ws.Cells.Item(27,8).Value = -2.925713

$ws.Cells.Item(28,1).Value = 'INE238AD6AF6'
$ws.Cells.Item(28,2).Value = 'Axis Bank Limited CD 07-Jan-2026'
$ws.Cells.Item(28,3).Value = 'quant Liquid Fund'
$ws.Cells.Item(28,4).Value = 0
$ws.Cells.Item(28,5).Value = 0
$ws.Cells.Item(28,6).Value = 2.908182
$ws.Cells.Item(28,7).Value = 0
$ws.Cells.Item(28,8).Value = -2.908182

Write-Output "Wrote 28 rows"